$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set final cell values (labels in column A, content in B/C) ---
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("B2").Value = "LOT2054"
$ws.Range("C2").Value = "LOT2054"
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Engenharia de Segurança do Trabalho e Biossegurança"
$ws.Range("C3").Value = " Engenharia de Segurança do Trabalho e Biossegurança"
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Safety Engineering and Biosafety"
$ws.Range("C4").Value = "Safety Engineering and Biosafety"
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "2"
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2018"
$ws.Range("C8").Value = "01/01/2018"
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EB-6"
$ws.Range("C9").Value = "EB-6"
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "8853480 - Tatiane da Franca Silva"
$ws.Range("C10").Value = "8853480 - Tatiane da Franca Silva"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Introduction to occupational health and safety;Technical knowledge and legal aspects in biosafety; Biosafety risk to the worker, the community, and the environment;"
$ws.Range("C14").Value = "Introduction to occupational health and safety;Technical knowledge and legal aspects in biosafety; Biosafety risk to the worker, the community, and the environment;"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Work safety management and strategies for prevention and risks control in the workplace. Regulatory norms;Assessment of biological, chemical and radiological hazard in biotechnology procedure;Handling and disposal of contaminated waste;Biosafety level criteria and Biosecurity regulations; Regulation of genetically modified organism and its products;Biosafety in laboratory animal handling;Principles of Bioethics;Case studies;"
$ws.Range("C16").Value = "Work safety management and strategies for prevention and risks control in the workplace. Regulatory norms;Assessment of biological, chemical and radiological hazard in biotechnology procedure;Handling and disposal of contaminated waste;Biosafety level criteria and Biosecurity regulations; Regulation of genetically modified organism and its products;Biosafety in laboratory animal handling;Principles of Bioethics;Case studies;"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8853480 - Tatiane da Franca Silva"
$ws.Range("C18").Value = "8853480 - Tatiane da Franca Silva"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Duas notas N1e N2 distribuídas ao longo do semestre. A composição das `"N`" fica critério dodocente."
$ws.Range("C19").Value = "Duas notas N1e N2 distribuídas ao longo do semestre. A composição das `"N`" fica critério dodocente."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "MF = (N1 + N2)/2"
$ws.Range("C20").Value = "MF = (N1 + N2)/2"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor doque 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor doque 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0."
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOT2040 -  Engenharia Genética  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOT2040 -  Engenharia Genética  (Requisito fraco)`n"

# --- Clear cells that no longer hold content in the final layout ---
$ws.Range("B17:C17").ClearContents()
$ws.Range("B22:C22").ClearContents()
$ws.Range("A23").ClearContents()

# --- Remove the now-unused trailing row so the sheet shrinks to A1:C23 ---
$ws.Rows.Item(24).Delete()

# --- Row heights to match final layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).EntireRow.AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).EntireRow.AutoFit()
$ws.Rows.Item(23).RowHeight = 30
